$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update metrics for year 2025 (row 6)
$ws.Range("C6").Value = 387
$ws.Range("E6").Value = 82
$ws.Range("G6").Value = 21.18863049095607
$ws.Range("H6").Value = 78.81136950904393
